# Append the latest fetched exchange-rate log entries.
# allData_sheet gets one new row per currency; each per-currency sheet
# (CNY/JPY/GBP/EUR/USD) gets one new row with that currency's latest
# rate + update timestamp.

$wb = $excel.ActiveWorkbook

$timestamp = "2018-03-27 14:13:00"

# currency -> latest rate (kept as plain text, like the existing log rows)
$rates = @{
    "CNY" = "0.1091"
    "JPY" = "1.8462"
    "GBP" = "0.0123"
    "EUR" = "0.0140"
    "USD" = "0.0175"
}
$order = @("CNY", "JPY", "GBP", "EUR", "USD")

function Set-TextValue($range, $value) {
    # Rate strings look numeric ("0.1091", "1.8462", ...); force them to be
    # stored as text (matching the source data) instead of being coerced to
    # a number, then drop back to the default style so no formatting sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# allData_sheet: header in row 1, existing data rows 2-6, append rows 7-11
$all = $wb.Worksheets.Item("allData_sheet")
$nextRow = 7
foreach ($code in $order) {
    $all.Cells.Item($nextRow, 1).Value = $code
    Set-TextValue $all.Cells.Item($nextRow, 2) $rates[$code]
    $all.Cells.Item($nextRow, 3).Value = $timestamp
    $nextRow++
}

# Per-currency sheets: header in row 1, existing data row 2, append row 3
foreach ($code in $order) {
    $sheet = $wb.Worksheets.Item($code)
    Set-TextValue $sheet.Cells.Item(3, 1) $rates[$code]
    $sheet.Cells.Item(3, 2).Value = $timestamp
}
